$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly-adjusted floating point timestamp in D13
$ws.Range("D13").Value = 46001.6425868287

# Add the three new comment rows
$ws.Range("A14").Value = "PROCESO DE REGISTRO Y CIERRE DE ATENCIÓN POR FALLECIMIENTO"
$ws.Range("B14").Value = "quiero revisar este proceso"
$ws.Range("C14").Value = "jair"
$ws.Range("D14").Value = 46014.35867870371
$ws.Range("E14").Value = 46014

$ws.Range("A15").Value = "PROCEDIMIENTOS DE PQRS "
$ws.Range("B15").Value = "quiero este"
$ws.Range("C15").Value = "jair"
$ws.Range("D15").Value = 46014.36158993056
$ws.Range("E15").Value = 46014

$ws.Range("A16").Value = "PROCEDIMIENTOS DE PQRS "
$ws.Range("B16").Value = "no da"
$ws.Range("C16").Value = "jair"
$ws.Range("D16").Value = 46014.36193435869
$ws.Range("E16").Value = 46014

# Apply the same number format as the other date cells (D/E columns) to the new rows
$ws.Range("D14:E16").NumberFormat = $ws.Range("D13").NumberFormat
